$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '50.990.16'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.46%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.901.92'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -2.86%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.00%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '369.13'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.86%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '103.20'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -6.68%  '
$ws.Range("E7").Value = '  -5.66%  '
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.589'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -7.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '36.90'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -6.44%  '
$ws.Range("E11").Value = '  +0.58%  '
$ws.Range("E12").Value = '  -5.55%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '18.36'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -6.52%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.357.62'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.38%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.36'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -6.19%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.894.84'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.69%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.942'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -5.17%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '50.906.85'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.39%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.28'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -5.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.24'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -5.82%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.92'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -8.09%  '
$ws.Range("E22").Value = '  -4.89%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '68.27'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '259.32'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -4.97%  '
$ws.Range("E25").Value = '  -4.35%  '
$ws.Range("E26").Value = '  +3.52%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.170'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -5.59%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.15%  '
$ws.Range("E29").Value = '  -5.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.32'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -7.62%  '
$ws.Range("E31").Value = '  -5.29%  '
$ws.Range("B32").Value = 'Cosmos'
$ws.Range("C32").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '9.89'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -6.41%  '
$ws.Range("B33").Value = 'RenderToken'
$ws.Range("C33").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '6.08'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.61%  '
$ws.Range("E34").Value = '  -2.14%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '35.01'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -8.45%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '51.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.83%  '
$ws.Range("E37").Value = '  +0.58%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.0420'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -6.36%  '
$ws.Range("E39").Value = '  -5.92%  '
$ws.Range("E40").Value = '  -4.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.93'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -7.67%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.84'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -9.48%  '
$ws.Range("E43").Value = '  -6.51%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '22.21'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -7.15%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '117.43'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.45%  '
$ws.Range("E46").Value = '  -3.57%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.042.90'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -5.21%  '
$ws.Range("E48").Value = '  -6.06%  '
$ws.Range("E49").Value = '  -9.30%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '3.198.12'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -2.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.241'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -3.36%  '
